$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename "Config" sheet to "Theme" and replace its content with the new
#    theme-color key/value table.
# ---------------------------------------------------------------------------
$wsTheme = $wb.Worksheets.Item("Config")
$wsTheme.Name = "Theme"
$wsTheme.Cells.Clear()

$wsTheme.Range("A1").Value = "BackgroundColor"
$wsTheme.Range("B1").Value = "TextColor"
$wsTheme.Range("C1").Value = "HeaderBackground"
$wsTheme.Range("D1").Value = "FooterBackground"

$wsTheme.Range("A2").Value = "#f0f4f8"
$wsTheme.Range("B2").Value = "#334155"
$wsTheme.Range("C2").Value = "#ffffff"
$wsTheme.Range("D2").Value = "#1e293b"

# ---------------------------------------------------------------------------
# 2) "Content" sheet: replace the Category/Title/Subtitle/Content/Copyright
#    table with a simpler Section/Text table.
# ---------------------------------------------------------------------------
$wsContent = $wb.Worksheets.Item("Content")
$wsContent.Cells.Clear()

$wsContent.Range("A1").Value = "Section"
$wsContent.Range("B1").Value = "Text"

$wsContent.Range("A2").Value = "HeroTitle"
$wsContent.Range("B2").Value = "See the World Clearly"

$wsContent.Range("A3").Value = "HeroSubtitle"
$wsContent.Range("B3").Value = "Premium eyewear and contact lenses for comfort and style."

$wsContent.Range("A4").Value = "FooterTagline"
$wsContent.Range("B4").Value = "Your Vision Partner"

# ---------------------------------------------------------------------------
# 3) "Navigation" sheet: only the A1 header text changes.
# ---------------------------------------------------------------------------
$wsNav = $wb.Worksheets.Item("Navigation")
$wsNav.Range("A1").Value = "MenuItem"

# ---------------------------------------------------------------------------
# 4) "Products" sheet: drop the Category column and refresh the product
#    rows; prices become real numbers instead of "$"-prefixed text.
# ---------------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Cells.Clear()

$wsProducts.Range("A1").Value = "ProductName"
$wsProducts.Range("B1").Value = "Description"
$wsProducts.Range("C1").Value = "Price"

$wsProducts.Range("A2").Value = "Classic Black Frames"
$wsProducts.Range("B2").Value = "Timeless acetate frames with durable hinges."
$wsProducts.Range("C2").Value = 129.99

$wsProducts.Range("A3").Value = "Air Light Titanium"
$wsProducts.Range("B3").Value = "Ultra-lightweight titanium for all-day comfort."
$wsProducts.Range("C3").Value = 189.99

$wsProducts.Range("A4").Value = "Blue Light Protection"
$wsProducts.Range("B4").Value = "Reduces eye strain from digital screens."
$wsProducts.Range("C4").Value = 159.99

# ---------------------------------------------------------------------------
# 5) "Services" sheet: header rename + refreshed rows; prices become real
#    numbers instead of "$"-prefixed text.
# ---------------------------------------------------------------------------
$wsServices = $wb.Worksheets.Item("Services")
$wsServices.Cells.Clear()

$wsServices.Range("A1").Value = "ServiceName"
$wsServices.Range("B1").Value = "Description"
$wsServices.Range("C1").Value = "Price"

$wsServices.Range("A2").Value = "Eye Exam"
$wsServices.Range("B2").Value = "Comprehensive vision testing and diagnosis."
$wsServices.Range("C2").Value = 49.99

$wsServices.Range("A3").Value = "Lens Fitting"
$wsServices.Range("B3").Value = "Personalized contact lens fitting session."
$wsServices.Range("C3").Value = 79.99

$wsServices.Range("A4").Value = "Frame Styling"
$wsServices.Range("B4").Value = "Expert guidance to choose the perfect frame."
$wsServices.Range("C4").Value = 29.99

# ---------------------------------------------------------------------------
# 6) "Contact" sheet: drop the Category column and refresh the actual
#    contact details.
# ---------------------------------------------------------------------------
$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Cells.Clear()

$wsContact.Range("A1").Value = "Address"
$wsContact.Range("B1").Value = "City"
$wsContact.Range("C1").Value = "Country"
$wsContact.Range("D1").Value = "Phone"
$wsContact.Range("E1").Value = "Email"
$wsContact.Range("F1").Value = "Hours"

$wsContact.Range("A2").Value = "123 Main Street, Kuala Lumpur"
$wsContact.Range("B2").Value = "Kuala Lumpur"
$wsContact.Range("C2").Value = "Malaysia"
$wsContact.Range("D2").Value = "(+60) 3-2212 3456"
$wsContact.Range("E2").Value = "info@cerminmataho.my"
$wsContact.Range("F2").Value = "Mon-Fri: 9AM–7PM, Sat-Sun: 10AM–6PM"
